# "Fix 'Rejeter avenant' issue"
# The sheet had two extra helper/summary rows (9 and 10) appended below the
# data table. Those rows are removed (shifting the footer/summary rows that
# used to live there up into rows 7-8) and the remaining data rows are
# updated with corrected values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two trailing rows (old rows 9 & 10) - this also shrinks the
# sheet's used range/dimension from A1:L10 down to A1:L8 automatically.
$ws.Rows("9:10").Delete()

# Make sure every cell we are about to fill keeps its numeric-looking text
# (amounts, dates, years, ...) stored as TEXT rather than being
# auto-converted to a number/date by Excel's input parser.
$ws.Range("A1:L8").NumberFormat = "@"

# --- Row 2 ---------------------------------------------------------------
$ws.Range("A2").Value = ""
$ws.Range("B2").Value = "K5443645"
$ws.Range("C2").Value = ""
$ws.Range("D2").Value = "KHADIJA LALA"
$ws.Range("E2").Value = "KHDEHOI35456"
$ws.Range("F2").Value = "AAAAAA"
$ws.Range("G2").Value = "LUC"
$ws.Range("H2").Value = "-"
$ws.Range("I2").Value = "100000.00"
$ws.Range("J2").Value = "13500.00"
$ws.Range("K2").Value = "86500.00"
$ws.Range("L2").Value = "TRS.15.2018"

# --- Row 3 ---------------------------------------------------------------
$ws.Range("A3").Value = ""
$ws.Range("B3").Value = "L3578354"
$ws.Range("C3").Value = ""
$ws.Range("D3").Value = "NABIL KAMAL"
$ws.Range("E3").Value = "PMLPL35434"
$ws.Range("F3").Value = "WWWWW"
$ws.Range("G3").Value = "LUC"
$ws.Range("H3").Value = "-"
$ws.Range("I3").Value = "70000.00"
$ws.Range("J3").Value = "9000.00"
$ws.Range("K3").Value = "61000.00"
$ws.Range("L3").Value = "TRS.15.2018"

# --- Row 4 ---------------------------------------------------------------
$ws.Range("A4").Value = ""
$ws.Range("B4").Value = "D524564"
$ws.Range("C4").Value = ""
$ws.Range("D4").Value = "SAMIRA TATA"
$ws.Range("E4").Value = "LKJOIFEJIOZ"
$ws.Range("F4").Value = "QW"
$ws.Range("G4").Value = "LUC"
$ws.Range("H4").Value = "-"
$ws.Range("I4").Value = "60000.00"
$ws.Range("J4").Value = "7500.00"
$ws.Range("K4").Value = "52500.00"
$ws.Range("L4").Value = "TRS.15.2018"

# --- Row 5 ---------------------------------------------------------------
$ws.Range("A5").Value = ""
$ws.Range("F5").Value = "QQQQQQQQQ"
$ws.Range("I5").Value = "50000.00"
$ws.Range("J5").Value = "6000.00"
$ws.Range("K5").Value = "44000.00"
$ws.Range("L5").Value = "TRS.15.2018"

# --- Row 6 (now fully blank) ----------------------------------------------
$ws.Range("A6:L6").Value = ""

# --- Row 7 (becomes the field-name header row) ----------------------------
$ws.Range("A7").Value = ""
$ws.Range("B7").Value = ""
$ws.Range("C7").Value = ""
$ws.Range("D7").Value = "identifiantFiscal"
$ws.Range("E7").Value = "exerciceFiscalDu"
$ws.Range("F7").Value = "exerciceFiscalAu"
$ws.Range("G7").Value = "annee"
$ws.Range("H7").Value = "totalMntBrutLoyer"
$ws.Range("I7").Value = "totalMntRetenueSource"
$ws.Range("J7").Value = "totalMntNetLoyer"
$ws.Range("K7").Value = ""
$ws.Range("L7").Value = ""

# --- Row 8 (becomes the totals row) ---------------------------------------
$ws.Range("D8").Value = "IF"
$ws.Range("E8").Value = "2022-01-01"
$ws.Range("F8").Value = "2022-12-31"
$ws.Range("G8").Value = "2022"
$ws.Range("H8").Value = "280000.00"
$ws.Range("I8").Value = "36000.00"
$ws.Range("J8").Value = "244000.00"
